$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.310.64"
$ws.Range("E2").Value = "  +1.05%  "
$ws.Range("D3").Value = "3.356.52"
$ws.Range("E3").Value = "  +0.69%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "584.17"
$ws.Range("E5").Value = "  +0.64%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "177.68"
$ws.Range("E6").Value = "  +1.05%  "
$ws.Range("E7").Value = "  -0.16%  "
$ws.Range("E8").Value = "  -0.08%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.185"
$ws.Range("E9").Value = "  +3.48%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.582"
$ws.Range("E10").Value = "  +0.95%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "48.11"
$ws.Range("E11").Value = "  +5.83%  "
$ws.Range("E12").Value = "  +1.78%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "688.73"
$ws.Range("E13").Value = "  +3.47%  "
$ws.Range("D14").Value = "3.896.68"
$ws.Range("E14").Value = "  +0.45%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.44"
$ws.Range("E15").Value = "  +0.35%  "
$ws.Range("D16").Value = "68.330.98"
$ws.Range("E16").Value = "  +0.86%  "
$ws.Range("E17").Value = "  +1.30%  "
$ws.Range("D18").Value = "3.354.39"
$ws.Range("E18").Value = "  +0.59%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "17.49"
$ws.Range("E19").Value = "  +0.63%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.25"
$ws.Range("E20").Value = "  +2.33%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.896"
$ws.Range("E21").Value = "  +0.76%  "
$ws.Range("E22").Value = "  +0.22%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "17.02"
$ws.Range("E23").Value = "  -0.32%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "100.55"
$ws.Range("E24").Value = "  +1.22%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.92"
$ws.Range("E25").Value = "  +1.59%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.71"
$ws.Range("E26").Value = "  +1.62%  "
$ws.Range("E27").Value = "  +2.32%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "33.05"
$ws.Range("E28").Value = "  -1.75%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.54"
$ws.Range("E29").Value = "  +1.29%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.99"
$ws.Range("E30").Value = "  -5.47%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "11.09"
$ws.Range("E31").Value = "  +1.07%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "555.64"
$ws.Range("E32").Value = "  -3.98%  "
$ws.Range("E33").Value = "  +0.70%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "58.06"
$ws.Range("E34").Value = "  +2.56%  "
$ws.Range("E35").Value = "  +0.09%  "
$ws.Range("D36").Value = "3.714.30"
$ws.Range("E36").Value = "  +0.44%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.37"
$ws.Range("E37").Value = "  -0.01%  "
$ws.Range("E38").Value = "  +4.80%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "34.83"
$ws.Range("E39").Value = "  +1.31%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.18"
$ws.Range("E40").Value = "  +2.10%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.62"
$ws.Range("E41").Value = "  -0.16%  "
$ws.Range("D42").Value = "0.0₃0676"
$ws.Range("E42").Value = "  +1.32%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.336"
$ws.Range("E43").Value = "  +0.65%  "
$ws.Range("E44").Value = "  -1.41%  "
$ws.Range("E45").Value = "  +1.62%  "
$ws.Range("E46").Value = "  +1.73%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.129"
$ws.Range("E47").Value = "  +0.50%  "
$ws.Range("E48").Value = "  -0.15%  "
$ws.Range("E49").Value = "  -0.65%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "131.91"
$ws.Range("E50").Value = "  +2.32%  "
$ws.Range("E51").Value = "  -1.86%  "
